$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.025.40'
$ws.Range("E2").Value = '  +0.56%  '
$ws.Range("D3").Value = '1.590.87'
$ws.Range("E3").Value = '  +0.32%  '
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.44%  '
$ws.Range("E6").Value = '  -0.21%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.479'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.47%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.247'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0612'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.63%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '17.98'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.61%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0809'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.42%  '
$ws.Range("D12").Value = '1.812.85'
$ws.Range("E12").Value = '  +0.40%  '
$ws.Range("D13").Value = '1.585.00'
$ws.Range("E13").Value = '  +0.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.99'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.95%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.510'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.27%  '
$ws.Range("D16").Value = '26.011.68'
$ws.Range("E16").Value = '  +0.59%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.15'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.64%  '
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '202.46'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.23'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.19'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.99'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.72%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.95'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +14.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.14'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.92%  '
$ws.Range("E26").Value = '  -0.21%  '
$ws.Range("E27").Value = '  -7.94%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.08'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.47'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.10%  '
$ws.Range("E30").Value = '  +0.27%  '
$ws.Range("E31").Value = '  +0.73%  '
$ws.Range("E32").Value = '  -0.04%  '
$ws.Range("E33").Value = '  -2.92%  '
$ws.Range("E34").Value = '  -1.19%  '
$ws.Range("E35").Value = '  -0.70%  '
$ws.Range("D36").Value = '1.126.35'
$ws.Range("E36").Value = '  +2.53%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0162'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +8.39%  '
$ws.Range("E38").Value = '  -0.18%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.32'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.32%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.786'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.37%  '
$ws.Range("E41").Value = '  -2.18%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.779'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.42%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.12'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.72%  '
$ws.Range("D44").Value = '1.723.23'
$ws.Range("E44").Value = '  +0.23%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.14'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.74%  '
$ws.Range("E46").Value = '  -1.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '53.44'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.37%  '
$ws.Range("E48").Value = '  -1.21%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.407'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.29%  '
$ws.Range("E50").Value = '  +0.09%  '
$ws.Range("D51").Value = '0.0₇0920'
$ws.Range("E51").Value = '  -17.79%  '
